$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.884.92"
$ws.Range("D3").Value = "3.280.70"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.53"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.81"
$ws.Range("E6").Value = "  -5.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.582"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").Value = "3.275.16"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.34"
$ws.Range("E12").Value = "  -3.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000268"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "686.88"
$ws.Range("E14").Value = "  +8.04%  "
$ws.Range("D15").Value = "3.807.47"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.25"
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("D17").Value = "66.970.94"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "3.276.92"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("E20").Value = "  -4.00%  "
$ws.Range("E21").Value = "  -3.01%  "
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.00"
$ws.Range("E23").Value = "  -4.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.14"
$ws.Range("E24").Value = "  +2.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.05"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("E27").Value = "  -2.70%  "
$ws.Range("E28").Value = "  -2.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.34"
$ws.Range("E29").Value = "  +7.23%  "
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "573.10"
$ws.Range("E32").Value = "  -3.71%  "
$ws.Range("D33").Value = "3.863.59"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.79"
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("E35").Value = "  -2.90%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.36"
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.30"
$ws.Range("E38").Value = "  -13.84%  "
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.38"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.58"
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("D43").Value = "0.0₃0666"
$ws.Range("E43").Value = "  -5.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.98"
$ws.Range("E44").Value = "  -5.60%  "
$ws.Range("E45").Value = "  -3.25%  "
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("E49").Value = "  +6.50%  "
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.00"
$ws.Range("E51").Value = "  -0.57%  "
